# "added from batch file" - the file that was uploaded on Sheet1 (column D,
# file2upload) is replaced with a newly generated report produced by the
# batch run.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

$newFile = "C:\Users\johnny.abouhaidar\Desktop\myfiles\katalon\Test Run Report 26.pdf"

# Every row that previously referenced the old uploaded file now points at
# the new one produced by the batch file.
$ws1.Range("D2").Value = $newFile
$ws1.Range("D3").Value = $newFile
$ws1.Range("D4").Value = $newFile
$ws1.Range("D6").Value = $newFile

# Sheet1 is where the user ended up after the edit.
$ws1.Activate()
$ws1.Range("D5").Select()
